# Weekly fruit/vegetable update: insert a new "Sandia" price record as row 16,
# shifting the existing rows 16-34 down to 17-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (pushes existing rows 16..34 to 17..35)
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 44650
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100112028
$ws.Cells.Item(16, 7).Value = "Sandia"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Segunda"
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 325
$ws.Cells.Item(16, 12).Value = 350
$ws.Cells.Item(16, 13).Value = 338
$ws.Cells.Item(16, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 338
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = "Hortaliza"
